$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the Range to hold a literal text value (never auto-coerced to a
    # number/date) the way typing an apostrophe-prefixed entry into Excel
    # does, then restore the "Normal" style so no stray quote-prefix
    # formatting is left behind on the cell.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "43.581.54"
$ws.Range("E2").Value = "  -0.84%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.232.70"
$ws.Range("E3").Value = "  -0.07%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
Set-TextValue "D5" "271.37"
$ws.Range("E5").Value = "  +4.50%  "

# Row 6 - Solana
Set-TextValue "D6" "93.34"
$ws.Range("E6").Value = "  +13.10%  "

# Row 7 - XRP
Set-TextValue "D7" "0.622"
$ws.Range("E7").Value = "  -0.67%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.618"
$ws.Range("E9").Value = "  +2.21%  "

# Row 10 - Avalanche
Set-TextValue "D10" "46.34"
$ws.Range("E10").Value = "  +4.36%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0920"
$ws.Range("E11").Value = "  -1.03%  "

# Row 12 - Polkadot
Set-TextValue "D12" "8.07"
$ws.Range("E12").Value = "  +14.16%  "

# Row 13 - TRON
Set-TextValue "D13" "0.104"
$ws.Range("E13").Value = "  +0.50%  "

# Row 14 - was Chainlink, now WrappedliquidstakedEther2.0 (rows 14/15 swapped)
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D14" "2.571.42"
$ws.Range("E14").Value = "  +0.18%  "

# Row 15 - was WrappedliquidstakedEther2.0, now Chainlink
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D15" "15.12"
$ws.Range("E15").Value = "  +3.19%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "2.254.84"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +1.16%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "43.560.22"
$ws.Range("E18").Value = "  -0.67%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -1.22%  "

# Row 20 - Uniswap
Set-TextValue "D20" "6.00"
$ws.Range("E20").Value = "  -0.84%  "

# Row 21 - Litecoin
Set-TextValue "D21" "70.39"
$ws.Range("E21").Value = "  -1.12%  "

# Row 22 - ImmutableX
Set-TextValue "D22" "2.33"
$ws.Range("E22").Value = "  -2.21%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "232.43"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "8.78"
$ws.Range("E24").Value = "  -5.61%  "

# Row 25 - Dai
Set-TextValue "D25" "0.999"
$ws.Range("E25").Value = "  -0.08%  "

# Row 26 - was PancakeSwap, now Cosmos (rows 26/27 swapped)
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D26" "11.26"
$ws.Range("E26").Value = "  +4.33%  "

# Row 27 - was Cosmos, now PancakeSwap
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D27" "2.49"
$ws.Range("E27").Value = "  +10.62%  "

# Row 28 - WEMIXToken
Set-TextValue "D28" "3.55"
$ws.Range("E28").Value = "  +5.59%  "

# Row 29 - InjectiveProtocol
Set-TextValue "D29" "39.52"
$ws.Range("E29").Value = "  -3.35%  "

# Row 30 - Toncoin
Set-TextValue "D30" "2.26"
$ws.Range("E30").Value = "  +2.07%  "

# Row 31 - Monero
Set-TextValue "D31" "173.02"
$ws.Range("E31").Value = "  +0.22%  "

# Row 32 - Hedera
Set-TextValue "D32" "0.0925"
$ws.Range("E32").Value = "  +3.30%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "20.75"
$ws.Range("E33").Value = "  +0.36%  "

# Row 34 - Filecoin
Set-TextValue "D34" "5.41"
$ws.Range("E34").Value = "  +0.79%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  +0.27%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  -4.16%  "

# Row 37 - VeChain
Set-TextValue "D37" "0.0350"
$ws.Range("E37").Value = "  -5.84%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  -5.59%  "

# Row 39 - NEARProtocol
$ws.Range("E39").Value = "  +14.56%  "

# Row 40 - Celestia
Set-TextValue "D40" "12.59"
$ws.Range("E40").Value = "  -3.58%  "

# Row 41 - LidoDAOToken
$ws.Range("E41").Value = "  +1.19%  "

# Row 42 - Algorand
$ws.Range("E42").Value = "  +6.28%  "

# Row 43 - MultiversX
Set-TextValue "D43" "62.77"
$ws.Range("E43").Value = "  -1.45%  "

# Row 44 - THORChain
Set-TextValue "D44" "5.37"
$ws.Range("E44").Value = "  -3.08%  "

# Row 45 - Cronos
Set-TextValue "D45" "0.0989"
$ws.Range("E45").Value = "  +0.05%  "

# Row 46 - FraxShare
Set-TextValue "D46" "8.40"
$ws.Range("E46").Value = "  -0.52%  "

# Row 47 - Aave
Set-TextValue "D47" "99.85"
$ws.Range("E47").Value = "  -4.32%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  +1.37%  "

# Row 49 - TrustWalletToken
$ws.Range("E49").Value = "  +2.08%  "

# Row 50 - WOONetwork
Set-TextValue "D50" "0.436"
$ws.Range("E50").Value = "  -2.45%  "

# Row 51 - Stacks
$ws.Range("E51").Value = "  -6.99%  "
